# GSC export refresh: the daily breadcrumb export rolled forward by one day.
# The oldest date (2025-10-11) drops off the front of the "Chart" sheet and
# a new day (2026-01-08) is appended at the bottom, with every date/value
# pair shifting up by one row to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = 90

# --- Shift column A (dates) and column C (validation counts) up by one row ---
for ($r = 2; $r -lt $lastRow; $r++) {
    $srcDate = $ws.Cells.Item($r + 1, 1)
    $dstDate = $ws.Cells.Item($r, 1)
    $srcDate.Copy($dstDate)

    $dstCount = $ws.Cells.Item($r, 3)
    $srcCount = $ws.Cells.Item($r + 1, 3)
    $dstCount.Value = $srcCount.Value()
}

# --- Write the newly appended day into the last row ---
# A helper cell is used so the new date is written as literal text (matching
# the existing shared-string date cells) instead of being auto-converted into
# a date serial number by the normal Value-assignment auto-detection.
$helper = $ws.Cells.Item(200, 10)
$helper.Formula = "=""2026-01-08"""
$helper.Copy()
$ws.Cells.Item($lastRow, 1).PasteSpecial(-4163)
$helper.Clear()

$ws.Cells.Item($lastRow, 3).Value = 27.0
